$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric,
# so Excel keeps them as text (matching the source workbook's inlineStr cells)
# instead of auto-converting to numbers.
$textCells = @("D4", "D5", "D6", "D8", "D9", "D12", "D16", "D17", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D30", "D31", "D34", "D35", "D36", "D37", "D39", "D40", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "46.031.39"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "2.378.04"
$ws.Range("E3").Value = "  +3.30%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "301.03"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "98.79"
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("E10").Value = "  -6.29%  "
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "7.15"
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "2.743.32"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").Value = "2.392.33"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").Value = "0.827"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "13.76"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "45.964.39"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -5.47%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "6.06"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "66.83"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "244.01"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").Value = "  -5.31%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "1.94"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").Value = "39.82"
$ws.Range("E27").Value = "  -10.99%  "
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "20.97"
$ws.Range("E30").Value = "  +4.19%  "
$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").Value = "3.74"
$ws.Range("E31").Value = "  +18.25%  "
$ws.Range("E32").Value = "  +6.88%  "
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("D34").Value = "147.23"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "0.0774"
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  +5.93%  "
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").Value = "15.12"
$ws.Range("E39").Value = "  -5.27%  "
$ws.Range("D40").Value = "3.89"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  -7.74%  "
$ws.Range("D43").Value = "1.932.78"
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "92.20"
$ws.Range("E45").Value = "  +4.72%  "
$ws.Range("E46").Value = "  -9.11%  "
$ws.Range("D47").Value = "8.53"
$ws.Range("E47").Value = "  +5.46%  "
$ws.Range("D48").Value = "0.186"
$ws.Range("E48").Value = "  -5.03%  "
$ws.Range("D49").Value = "99.09"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").Value = "2.611.73"
$ws.Range("E50").Value = "  +3.28%  "
$ws.Range("D51").Value = "69.05"
$ws.Range("E51").Value = "  -7.24%  "
